$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.118.17'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.420.47'
$ws.Range('E3').Value = '  -3.48%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '488.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +18.90%  '
$ws.Range('D9').Value = '2.444.42'
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.20%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '2.836.67'
$ws.Range('E14').Value = '  -3.74%  '
$ws.Range('D15').Value = '57.130.55'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '2.442.05'
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '324.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('D28').Value = '2.525.66'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.59%  '
$ws.Range('D30').Value = '0.0₃0785'
$ws.Range('E30').Value = '  -5.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.819'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -10.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '285.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.80%  '
$ws.Range('E40').Value = '  +6.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '34.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.599'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('E46').Value = '  -5.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0227'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.12%  '
$ws.Range('D50').Value = '1.912.68'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.56%  '
